{"js": "// Replace the 25 division-expression text runs in the practice table with\n// their new values, matching the committed OOXML diff. Each \"old\" value is\n// unique in the document, so a literal, case-sensitive search for it\n// unambiguously locates the single run to update.\nconst replacements = [\n  [\"58\u00f74=\", \"42\u00f73=\"],\n  [\"21\u00f76=\", \"87\u00f76=\"],\n  [\"93\u00f78=\", \"80\u00f78=\"],\n  [\"20\u00f76=\", \"37\u00f76=\"],\n  [\"83\u00f78=\", \"63\u00f72=\"],\n  [\"17\u00f74=\", \"78\u00f79=\"],\n  [\"59\u00f74=\", \"74\u00f75=\"],\n  [\"27\u00f74=\", \"97\u00f78=\"],\n  [\"90\u00f76=\", \"38\u00f79=\"],\n  [\"61\u00f77=\", \"20\u00f79=\"],\n  [\"34\u00f72=\", \"94\u00f75=\"],\n  [\"79\u00f78=\", \"69\u00f75=\"],\n  [\"11\u00f72=\", \"81\u00f77=\"],\n  [\"52\u00f76=\", \"12\u00f79=\"],\n  [\"58\u00f77=\", \"12\u00f76=\"],\n  [\"96\u00f73=\", \"48\u00f79=\"],\n  [\"78\u00f73=\", \"21\u00f76=\"],\n  [\"41\u00f75=\", \"75\u00f76=\"],\n  [\"28\u00f77=\", \"90\u00f77=\"],\n  [\"51\u00f75=\", \"81\u00f78=\"],\n  [\"15\u00f74=\", \"65\u00f74=\"],\n  [\"73\u00f73=\", \"12\u00f76=\"],\n  [\"67\u00f73=\", \"30\u00f72=\"],\n  [\"19\u00f79=\", \"17\u00f77=\"],\n  [\"52\u00f73=\", \"22\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find expected text \"${oldText}\" to replace.`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the 25 division-expression text runs in the practice table with\n# their new values, matching the committed OOXML diff. Each \"old\" value is\n# unique in the document, so Find/Replace for the literal text unambiguously\n# targets the single run to update.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"58\u00f74=\", \"42\u00f73=\"),\n    @(\"21\u00f76=\", \"87\u00f76=\"),\n    @(\"93\u00f78=\", \"80\u00f78=\"),\n    @(\"20\u00f76=\", \"37\u00f76=\"),\n    @(\"83\u00f78=\", \"63\u00f72=\"),\n    @(\"17\u00f74=\", \"78\u00f79=\"),\n    @(\"59\u00f74=\", \"74\u00f75=\"),\n    @(\"27\u00f74=\", \"97\u00f78=\"),\n    @(\"90\u00f76=\", \"38\u00f79=\"),\n    @(\"61\u00f77=\", \"20\u00f79=\"),\n    @(\"34\u00f72=\", \"94\u00f75=\"),\n    @(\"79\u00f78=\", \"69\u00f75=\"),\n    @(\"11\u00f72=\", \"81\u00f77=\"),\n    @(\"52\u00f76=\", \"12\u00f79=\"),\n    @(\"58\u00f77=\", \"12\u00f76=\"),\n    @(\"96\u00f73=\", \"48\u00f79=\"),\n    @(\"78\u00f73=\", \"21\u00f76=\"),\n    @(\"41\u00f75=\", \"75\u00f76=\"),\n    @(\"28\u00f77=\", \"90\u00f77=\"),\n    @(\"51\u00f75=\", \"81\u00f78=\"),\n    @(\"15\u00f74=\", \"65\u00f74=\"),\n    @(\"73\u00f73=\", \"12\u00f76=\"),\n    @(\"67\u00f73=\", \"30\u00f72=\"),\n    @(\"19\u00f79=\", \"17\u00f77=\"),\n    @(\"52\u00f73=\", \"22\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    [void]$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
